# Updated symbol list on Thu Jan 26 03:15:13 UTC 2023 with GitHub Actions
# Applies the per-row Price (D) / Volume(1h) (E) / Hora (G) updates from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{R=2; D="307.39"; E="2.68%"; G="3"},
    @{R=3; D="36.15"; E="3.06%"; G="3"},
    @{R=4; D="5.097"; E="2.46%"; G="3"},
    @{R=5; D="0.08145"; E="3.59%"; G="3"},
    @{R=6; D="1.925"; E="1.18%"; G="3"},
    @{R=7; D="4.191"; E="4.08%"; G="3"},
    @{R=8; D="7.770"; E="0.58%"; G="3"},
    @{R=9; D="0.9290"; E="0.56%"; G="3"},
    @{R=10; E="26.14%"; G="3"},
    @{R=11; D="0.1928"; E="6.33%"; G="3"},
    @{R=12; D="0.09306"; E="1.06%"; G="3"},
    @{R=13; D="0.03569"; E="0.30%"; G="3"},
    @{R=14; D="0.09856"; E="-0.29%"; G="3"},
    @{R=15; D="0.001416"; E="0.85%"; G="3"},
    @{R=16; D="0.005767"; E="-0.15%"; G="3"},
    @{R=17; D="3.558"; E="2.08%"; G="3"},
    @{R=18; D="2.975"; E="2.30%"; G="3"},
    @{R=19; D="0.3438"; E="-0.10%"; G="3"},
    @{R=20; D="0.1303"; E="-0.45%"; G="3"},
    @{R=21; E="-2.81%"; G="3"},
    @{R=22; D="0.2406"; E="0.19%"; G="3"},
    @{R=23; D="0.04517"; E="-0.24%"; G="3"},
    @{R=24; D="0.001213"; E="-0.10%"; G="3"},
    @{R=25; D="0.004889"; E="6.71%"; G="3"},
    @{R=26; D="0.0001239"; E="-0.86%"; G="3"},
    @{R=27; G="3"},
    @{R=28; G="3"},
    @{R=29; G="3"},
    @{R=30; G="3"},
    @{R=31; G="3"},
    @{R=32; G="3"},
    @{R=33; G="3"},
    @{R=34; G="3"},
    @{R=35; G="3"},
    @{R=36; G="3"},
    @{R=37; G="3"},
    @{R=38; G="3"},
    @{R=39; D="0.02020"; E="7.91%"; G="3"},
    @{R=40; D="0.04942"; E="5.68%"; G="3"},
    @{R=41; D="0.01108"; E="15.81%"; G="3"},
    @{R=42; D="0.007707"; E="1.60%"; G="3"},
    @{R=43; D="0.1381"; E="4.60%"; G="3"},
    @{R=44; D="0.002098"; E="-1.00%"; G="3"},
    @{R=45; D="0.01050"; E="-3.25%"; G="3"},
    @{R=46; D="0.00006436"; E="7.18%"; G="3"},
    @{R=47; D="0.00000000749"; E="-0.19%"; G="3"},
    @{R=48; G="3"},
    @{R=49; D="0.001189"; E="-8.81%"; G="3"},
    @{R=50; D="0.00002098"; E="-0.19%"; G="3"},
    @{R=51; D="0.0001998"; E="-0.19%"; G="3"}
)

foreach ($item in $updates) {
    $r = $item.R
    if ($item.ContainsKey("D")) {
        $ws.Range("D$r").Value = "'" + $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E$r").Value = "'" + $item.E
    }
    if ($item.ContainsKey("G")) {
        $ws.Range("G$r").Value = "'" + $item.G
    }
}

Write-Output "Applied $($updates.Count) row updates"
